$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely, shifting B:F left to A:E
$ws.Range("A1:A5").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
